# Automatic update of files.
# Rewrites the data (non-location) fields for rows 25-33 on the active
# worksheet so that each row's species/observation data is shifted to a
# different row per the source-control diff, while row 29 stays the same.
#
# Mapping of "new row" -> "row whose data it now holds" (source row in the
# original workbook):
#   25 <- 26
#   26 <- 32
#   27 <- 25
#   28 <- 33
#   29 <- 29 (unchanged)
#   30 <- 28
#   31 <- 27
#   32 <- 30
#   33 <- 31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together with the row's data: A,B,E,F,G,H,Q,R,Z,AB
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

$mapping = @{
    25 = 26
    26 = 32
    27 = 25
    28 = 33
    29 = 29
    30 = 28
    31 = 27
    32 = 30
    33 = 31
}

# Snapshot the original values for every relevant column/row BEFORE any
# writes happen, so that later writes don't clobber data we still need to
# read from (important because row 29 stays in place and several source
# rows are read after other rows have already been overwritten).
$snapshot = @{}
foreach ($col in $cols) {
    for ($r = 25; $r -le 33; $r++) {
        $addr = "$col$r"
        $snapshot[$addr] = $ws.Range($addr).Value()
    }
}

foreach ($destRow in 25..33) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value = $snapshot[$srcAddr]
    }
}
